$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 2535.2222
$ws.Range("I4").Value = 2535.2222
$ws.Range("K4").Value = 2535.2222
$ws.Range("M4").Value = -2421.2222
$ws.Range("H13").Value = 1040
$ws.Range("I13").Value = 2120
$ws.Range("K13").Value = 2120
$ws.Range("M13").Value = -1951
$ws.Range("H33").Value = 150.6923
$ws.Range("I33").Value = 150.6923
$ws.Range("K33").Value = 150.6923
$ws.Range("M33").Value = 78.30770000000001
$ws.Range("H40").Value = 6381
$ws.Range("I40").Value = 3365.1667
$ws.Range("K40").Value = 3365.1667
$ws.Range("M40").Value = -3190.1667
$ws.Range("H41").Value = 2654.7778
$ws.Range("I41").Value = 26
$ws.Range("K41").Value = 26
$ws.Range("M41").Value = 414
$ws.Range("H96").Value = 2279.4443
$ws.Range("I96").Value = 2703.75
$ws.Range("K96").Value = 8111.25
$ws.Range("M96").Value = -6738.25
$ws.Range("H107").Value = 959.4737
$ws.Range("I107").Value = 612.25
$ws.Range("J107").Value = 1554.7142
$ws.Range("K107").Value = 612.25
$ws.Range("L107").Value = 1554.7142
$ws.Range("M107").Value = 1307.75
$ws.Range("N107").Value = -5394.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6750
$ws.Range("I2").Value = 6500
$ws.Range("K2").Value = 6500
$ws.Range("M2").Value = -6387
$ws.Range("H33").Value = 4000
$ws.Range("I33").Value = 4000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 4000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -3671
$ws.Range("N33").ClearContents()
$ws.Range("H44").Value = 40000
$ws.Range("I44").Value = 45000
$ws.Range("J44").Value = 35000
$ws.Range("K44").Value = 45000
$ws.Range("L44").Value = 35000
$ws.Range("M44").Value = -44512
$ws.Range("N44").Value = -35976
$ws.Range("H63").Value = 1557.8572
$ws.Range("I63").Value = 1557.8572
$ws.Range("K63").Value = 1557.8572
$ws.Range("M63").Value = -871.8571999999999
$ws.Range("H66").Value = 1557.8572
$ws.Range("I66").Value = 1557.8572
$ws.Range("K66").Value = 7789.286
$ws.Range("M66").Value = -4357.286
$ws.Range("H116").Value = 6750
$ws.Range("I116").Value = 6500
$ws.Range("K116").Value = 6500
$ws.Range("M116").Value = -4206
$ws.Range("H135").Value = 39999.5
$ws.Range("J135").Value = 39999.5
$ws.Range("L135").Value = 39999.5
$ws.Range("N135").Value = -50139.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6750
$ws.Range("I3").Value = 6500
$ws.Range("K3").Value = 6500
$ws.Range("M3").Value = -6386
$ws.Range("H99").Value = 6559.8
$ws.Range("I99").Value = 6559.8
$ws.Range("K99").Value = 6559.8
$ws.Range("M99").Value = -5061.8
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 313.16666
$ws.Range("I10").Value = 313.16666
$ws.Range("K10").Value = 313.16666
$ws.Range("M10").Value = -174.16666
$ws.Range("H35").Value = 1114
$ws.Range("I35").Value = 1317.5
$ws.Range("J35").Value = 300
$ws.Range("K35").Value = 1317.5
$ws.Range("L35").Value = 300
$ws.Range("M35").Value = -1023.5
$ws.Range("N35").Value = -888
$ws.Range("H62").Value = 8441.154
$ws.Range("I62").Value = 9224
$ws.Range("K62").Value = 9224
$ws.Range("M62").Value = -8600
$ws.Range("H65").Value = 8441.154
$ws.Range("I65").Value = 9224
$ws.Range("K65").Value = 46120
$ws.Range("M65").Value = -43000
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2548.7856
$ws.Range("I5").Value = 1621
$ws.Range("K5").Value = 4863
$ws.Range("M5").Value = -4751
$ws.Range("H63").Value = 728
$ws.Range("I63").Value = 728
$ws.Range("K63").Value = 2184
$ws.Range("M63").Value = -1435
$ws.Range("H66").Value = 728
$ws.Range("I66").Value = 728
$ws.Range("K66").Value = 6552
$ws.Range("M66").Value = -2808
$ws.Range("H135").Value = 2548.7856
$ws.Range("I135").Value = 1621
$ws.Range("K135").Value = 14589
$ws.Range("M135").Value = -12054

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3038.5
$ws.Range("I102").Value = 2719.818
$ws.Range("K102").Value = 2719.818
$ws.Range("M102").Value = -1097.818
$ws.Range("H107").Value = 786.25
$ws.Range("I107").Value = 681.8333
$ws.Range("J107").Value = 1099.5
$ws.Range("K107").Value = 681.8333
$ws.Range("L107").Value = 1099.5
$ws.Range("M107").Value = 1238.1667
$ws.Range("N107").Value = -4939.5
$ws.Range("H122").Value = 5952.4165
$ws.Range("I122").Value = 6129.1816
$ws.Range("K122").Value = 18387.5448
$ws.Range("M122").Value = -15937.5448
$ws.Range("H132").Value = 3050.353
$ws.Range("I132").Value = 3337.0667
$ws.Range("J132").Value = 900
$ws.Range("K132").Value = 10011.2001
$ws.Range("L132").Value = 2700
$ws.Range("M132").Value = -7481.2001
$ws.Range("N132").Value = -7760

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 500
$ws.Range("I32").Value = 500
$ws.Range("K32").Value = 500
$ws.Range("M32").Value = -183
$ws.Range("H82").Value = 1529.5333
$ws.Range("I82").Value = 1559.6
$ws.Range("J82").Value = 1469.4
$ws.Range("K82").Value = 1559.6
$ws.Range("L82").Value = 1469.4
$ws.Range("M82").Value = -1198.6
$ws.Range("N82").Value = -2191.4
$ws.Range("H85").Value = 1529.5333
$ws.Range("I85").Value = 1559.6
$ws.Range("J85").Value = 1469.4
$ws.Range("K85").Value = 1559.6
$ws.Range("L85").Value = 1469.4
$ws.Range("M85").Value = -311.5999999999999
$ws.Range("N85").Value = -3965.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H81").Value = 2487.7
$ws.Range("I81").Value = 1982.7142
$ws.Range("K81").Value = 3965.4284
$ws.Range("M81").Value = -2904.4284
$ws.Range("H84").Value = 2487.7
$ws.Range("I84").Value = 1982.7142
$ws.Range("K84").Value = 19827.142
$ws.Range("M84").Value = -14523.142
$ws.Range("H107").Value = 621.3077
$ws.Range("I107").Value = 527.8
$ws.Range("K107").Value = 1583.4
$ws.Range("M107").Value = 336.6000000000001
$ws.Range("H136").Value = 2438.8
$ws.Range("I136").Value = 2519.7856
$ws.Range("K136").Value = 7559.3568
$ws.Range("M136").Value = -5009.3568
